$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset for rows 2..40 of columns A (index), B (wavelength), C (value)
$data = @(
    @(0, 417.7, 386.146387001185),
    @(1, 452.1, 378.795218929251),
    @(2, 486.8, 372.4852733461015),
    @(3, 521.5, 365.9174556086553),
    @(4, 556, 361.2590276644845),
    @(5, 563.5, 360.0819999357574),
    @(6, 588.4000000000001, 357.6899955828328),
    @(7, 614, 354.5517546411636),
    @(8, 639, 351.9624390873055),
    @(9, 664.6, 349.0976216976582),
    @(10, 689.5, 347.5243402529804),
    @(11, 715, 345.2350332250331),
    @(12, 740, 343.3188470134407),
    @(13, 765, 341.5650832514494),
    @(14, 790.5, 340.2582769739259),
    @(15, 816, 338.746861401007),
    @(16, 841.3000000000001, 337.01199303382),
    @(17, 866, 335.3928283401642),
    @(18, 891.5999999999999, 334.4406919724152),
    @(19, 917, 333.6484009438534),
    @(20, 942.4, 332.9635450431967),
    @(21, 967.3000000000001, 332.0792992030919),
    @(22, 992.7, 331.8855227337769),
    @(23, 1018, 331.0146621335886),
    @(24, 1060, 329.4586268876021),
    @(25, 1098, 328.6091423968518),
    @(26, 1138, 328.2417248821057),
    @(27, 1179, 327.4590579752656),
    @(28, 1222, 326.8927466910654),
    @(29, 1266, 326.1544314673292),
    @(30, 1312, 325.6030551689326),
    @(31, 1358, 324.8307467522693),
    @(32, 1408, 323.3997282622303),
    @(33, 1459, 322.2521987418236),
    @(34, 1512, 320.5516136163267),
    @(35, 1566, 318.8058571090835),
    @(36, 1623, 317.1814621512059),
    @(37, 1682, 316.8743183339373),
    @(38, 1743, 315.8021888031159)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Two brand-new rows (39 and 40) were added at the bottom; give column A of
# those rows the same formatting (bold/border/center) used by the rest of
# column A by copying the format from the last previously-existing row.
$ws.Cells.Item(38, 1).Copy() | Out-Null
$ws.Range("A39:A40").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
